$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.166.32'
$ws.Range("E2").Value = '  -0.80%  '
$ws.Range("D3").Value = '3.533.78'
$ws.Range("E3").Value = '  +2.51%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.74'
$ws.Range("E5").Value = '  +1.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.87'
$ws.Range("E6").Value = '  +0.27%  '
$ws.Range("D7").Value = '3.532.14'
$ws.Range("E7").Value = '  +2.48%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.489'
$ws.Range("E9").Value = '  -2.75%  '
$ws.Range("E10").Value = '  +1.49%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.91'
$ws.Range("E11").Value = '  -6.05%  '
$ws.Range("E12").Value = '  +2.53%  '
$ws.Range("D13").Value = '4.136.02'
$ws.Range("E13").Value = '  +2.80%  '
$ws.Range("E14").Value = '  +1.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.09'
$ws.Range("E15").Value = '  +1.81%  '
$ws.Range("D16").Value = '3.538.13'
$ws.Range("E16").Value = '  +2.96%  '
$ws.Range("D18").Value = '65.267.36'
$ws.Range("E18").Value = '  -0.51%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.30'
$ws.Range("E19").Value = '  +4.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.93'
$ws.Range("E20").Value = '  +0.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.30'
$ws.Range("E21").Value = '  +3.63%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '392.19'
$ws.Range("E22").Value = '  -0.74%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.573'
$ws.Range("E23").Value = '  +2.96%  '
$ws.Range("D24").Value = '3.675.15'
$ws.Range("E24").Value = '  +2.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '73.85'
$ws.Range("E25").Value = '  +0.50%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("E27").Value = '  +6.73%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.70'
$ws.Range("E28").Value = '  +6.83%  '
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("E30").Value = '  +1.56%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.21'
$ws.Range("E31").Value = '  -1.30%  '
$ws.Range("D32").Value = '3.547.28'
$ws.Range("E32").Value = '  +2.78%  '
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.76'
$ws.Range("E34").Value = '  +3.03%  '
$ws.Range("E35").Value = '  -0.87%  '
$ws.Range("E36").Value = '  +6.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.95'
$ws.Range("E37").Value = '  -0.28%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '168.69'
$ws.Range("E38").Value = '  -2.32%  '
$ws.Range("E39").Value = '  +3.30%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.95'
$ws.Range("E40").Value = '  +2.63%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0803'
$ws.Range("E41").Value = '  +4.26%  '
$ws.Range("E42").Value = '  -0.37%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.25'
$ws.Range("E43").Value = '  +13.84%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.83'
$ws.Range("E44").Value = '  -2.24%  '
$ws.Range("E45").Value = '  +0.14%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.42'
$ws.Range("E46").Value = '  -0.40%  '
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.67'
$ws.Range("E47").Value = '  +2.49%  '
$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.19'
$ws.Range("E48").Value = '  +5.80%  '
$ws.Range("D49").Value = '2.415.23'
$ws.Range("E49").Value = '  +8.97%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.79'
$ws.Range("E50").Value = '  +2.62%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '301.47'
$ws.Range("E51").Value = '  +6.45%  '
